# Generate Report for Handback
# Adds a second handed-back file (c7aebfb3-e028-4625-a464-c4a202a9e3ed.md) to the
# Overview / zh-cn / de-de report sheets, and refreshes the timestamps / xlf file
# names for the first file (renamed from d0ca1dd6-2ade-4e81-a460-1e79887cd4d5.md
# to 764aa2b8-bf11-4cf0-8544-6c5e104ae578.md).

$wb = $excel.ActiveWorkbook

$newGuid1 = "764aa2b8-bf11-4cf0-8544-6c5e104ae578"
$newGuid2 = "c7aebfb3-e028-4625-a464-c4a202a9e3ed"

$overviewDate = "2016-08-22 17:03:57"

$zhHash1 = "936d422752ada39154efa47a1f3bdfe180c42165"
$zhHash2 = "b4b684eb810b3565fcd59ccf816335e12673bde2"

$zhStart1 = "2016-08-22 17:03:53"
$zhEnd1   = "2016-08-22 17:04:20"
$deEnd1   = "2016-08-22 17:04:27"

function Set-Hyperlink($ws, $rangeRef, $target, $displayText) {
    $rng = $ws.Range($rangeRef)
    # NB: calling .Hyperlinks.Delete() on a range that has *no* hyperlink
    # wipes out every hyperlink on the sheet, so only call it when needed.
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $target, "", "", $displayText) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# refresh row 2 (existing file) with its new name + generate date
$wsOv.Range("A2").Value2 = "$newGuid1.md"
Set-Hyperlink $wsOv "B2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid1.md" "e2e\$newGuid1.md"
$wsOv.Range("G2").Value2 = $overviewDate

# add row 3 (new file)
$wsOv.Range("A3").Value2 = "$newGuid2.md"
$wsOv.Range("A3").Style = $wsOv.Range("A2").Style
Set-Hyperlink $wsOv "B3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid2.md" "e2e\$newGuid2.md"
$wsOv.Range("C3").Value2 = $wsOv.Range("C2").Value2
$wsOv.Range("E3").Value2 = $wsOv.Range("E2").Value2
$wsOv.Range("F3").Value2 = $wsOv.Range("F2").Value2
$wsOv.Range("G3").Value2 = $overviewDate
$wsOv.Range("G3").NumberFormat = $wsOv.Range("G2").NumberFormat

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf1 = "$newGuid1.$zhHash1.zh-cn.xlf"
$zhXlf2 = "$newGuid2.$zhHash2.zh-cn.xlf"

# refresh row 2 (existing file)
Set-Hyperlink $wsZh "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid1.md" "$newGuid1.md"
$wsZh.Range("G2").Value2 = $zhXlf1
$wsZh.Range("H2").Value2 = $zhStart1
Set-Hyperlink $wsZh "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid1.md" "$newGuid1.md"
$wsZh.Range("J2").Value2 = $zhXlf1
$wsZh.Range("K2").Value2 = $zhEnd1

# add row 3 (new file)
Set-Hyperlink $wsZh "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid2.md" "$newGuid2.md"
$wsZh.Range("A3").Style = $wsZh.Range("A2").Style
$wsZh.Range("B3").Value2 = $wsZh.Range("B2").Value2
$wsZh.Range("C3").Value2 = $wsZh.Range("C2").Value2
$wsZh.Range("D3").Value2 = $wsZh.Range("D2").Value2
$wsZh.Range("E3").Value2 = $wsZh.Range("E2").Value2
$wsZh.Range("F3").Value2 = "True"
$wsZh.Range("G3").Value2 = $zhXlf2
$wsZh.Range("H3").Value2 = $zhStart1
$wsZh.Range("H3").NumberFormat = $wsZh.Range("H2").NumberFormat
Set-Hyperlink $wsZh "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid2.md" "$newGuid2.md"
$wsZh.Range("I3").Style = $wsZh.Range("I2").Style
$wsZh.Range("J3").Value2 = $zhXlf2
$wsZh.Range("K3").Value2 = $zhEnd1
$wsZh.Range("K3").NumberFormat = $wsZh.Range("K2").NumberFormat
$wsZh.Range("L3").Value2 = $wsZh.Range("L2").Value2
$wsZh.Range("M3").Value2 = $wsZh.Range("M2").Value2
$wsZh.Range("N3").Value2 = $wsZh.Range("N2").Value2
$wsZh.Range("O3").Value2 = $wsZh.Range("O2").Value2
$wsZh.Range("P3").Value2 = $wsZh.Range("P2").Value2

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf1 = "$newGuid1.$zhHash1.de-de.xlf"
$deXlf2 = "$newGuid2.$zhHash2.de-de.xlf"

# refresh row 2 (existing file)
Set-Hyperlink $wsDe "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid1.md" "$newGuid1.md"
$wsDe.Range("G2").Value2 = $deXlf1
$wsDe.Range("H2").Value2 = $overviewDate
Set-Hyperlink $wsDe "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid1.md" "$newGuid1.md"
$wsDe.Range("J2").Value2 = $deXlf1
$wsDe.Range("K2").Value2 = $deEnd1

# add row 3 (new file)
Set-Hyperlink $wsDe "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid2.md" "$newGuid2.md"
$wsDe.Range("A3").Style = $wsDe.Range("A2").Style
$wsDe.Range("B3").Value2 = $wsDe.Range("B2").Value2
$wsDe.Range("C3").Value2 = $wsDe.Range("C2").Value2
$wsDe.Range("D3").Value2 = $wsDe.Range("D2").Value2
$wsDe.Range("E3").Value2 = $wsDe.Range("E2").Value2
$wsDe.Range("F3").Value2 = "True"
$wsDe.Range("G3").Value2 = $deXlf2
$wsDe.Range("H3").Value2 = $overviewDate
$wsDe.Range("H3").NumberFormat = $wsDe.Range("H2").NumberFormat
Set-Hyperlink $wsDe "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$newGuid2.md" "$newGuid2.md"
$wsDe.Range("I3").Style = $wsDe.Range("I2").Style
$wsDe.Range("J3").Value2 = $deXlf2
$wsDe.Range("K3").Value2 = $deEnd1
$wsDe.Range("K3").NumberFormat = $wsDe.Range("K2").NumberFormat
$wsDe.Range("L3").Value2 = $wsDe.Range("L2").Value2
$wsDe.Range("M3").Value2 = $wsDe.Range("M2").Value2
$wsDe.Range("N3").Value2 = $wsDe.Range("N2").Value2
$wsDe.Range("O3").Value2 = $wsDe.Range("O2").Value2
$wsDe.Range("P3").Value2 = $wsDe.Range("P2").Value2

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
